# Update the initial stock counts on the Medication sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PARACETAMOL initial stock: 19 -> 18
$ws.Range("B2").Value = 18

# AMOXICILLIN initial stock: 40 -> 39
$ws.Range("B4").Value = 39
